# Update countries & provincias Spain
# Applies the data refresh described in the commit:
#  - Reorders Martinica / Islas Feroe ahead of Madagascar in the country list
#    (rows 144-146 now show Martinica, Islas Feroe, Madagascar respectively,
#    each carrying the stats that belong to that country).
#  - Refreshes the daily COVID figures for Estados Unidos, India, Israel,
#    Uzbekistan, Ruanda, Martinica, Islas Feroe, Madagascar and Gibraltar.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 - Estados Unidos
$ws.Cells.Item(4,2).Value = 1378090
$ws.Cells.Item(4,3).Value = 10452
$ws.Cells.Item(4,4).Value = 259092
$ws.Cells.Item(4,5).Value = 1037773
$ws.Cells.Item(4,6).Value = 16551
$ws.Cells.Item(4,7).Value = 438
$ws.Cells.Item(4,8).Value = 81225

# Row 15 - India
$ws.Cells.Item(15,2).Value = 70768
$ws.Cells.Item(15,3).Value = 3607
$ws.Cells.Item(15,5).Value = 45925

# Row 33 - Israel
$ws.Cells.Item(33,6).Value = 58

# Row 75 - Uzbekistan
$ws.Cells.Item(75,2).Value = 2486
$ws.Cells.Item(75,3).Value = 68
$ws.Cells.Item(75,4).Value = 1988
$ws.Cells.Item(75,5).Value = 488

# Row 137 - Ruanda
$ws.Cells.Item(137,2).Value = 285
$ws.Cells.Item(137,3).Value = 1
$ws.Cells.Item(137,4).Value = 150
$ws.Cells.Item(137,5).Value = 135

# Rows 144-146: Martinica / Islas Feroe now sort ahead of Madagascar.
# Row 144 becomes Martinica
$ws.Cells.Item(144,1).Value = "Martinica"
$ws.Cells.Item(144,2).Value = 187
$ws.Cells.Item(144,3).Value = 1
$ws.Cells.Item(144,4).Value = 83
$ws.Cells.Item(144,5).Value = 90
$ws.Cells.Item(144,6).Value = 4
$ws.Cells.Item(144,7).Value = 0
$ws.Cells.Item(144,8).Value = 14

# Row 145 becomes Islas Feroe
$ws.Cells.Item(145,1).Value = "Islas Feroe"
$ws.Cells.Item(145,2).Value = 187
$ws.Cells.Item(145,3).Value = 0
$ws.Cells.Item(145,4).Value = 187
$ws.Cells.Item(145,5).Value = 0
$ws.Cells.Item(145,6).Value = 0
$ws.Cells.Item(145,7).Value = 0
$ws.Cells.Item(145,8).Value = 0

# Row 146 becomes Madagascar (with refreshed stats)
$ws.Cells.Item(146,1).Value = "Madagascar"
$ws.Cells.Item(146,2).Value = 186
$ws.Cells.Item(146,3).Value = 0
$ws.Cells.Item(146,4).Value = 101
$ws.Cells.Item(146,5).Value = 85
$ws.Cells.Item(146,6).Value = 1
$ws.Cells.Item(146,7).Value = 0
$ws.Cells.Item(146,8).Value = 0

# Row 153 - Gibraltar
$ws.Cells.Item(153,2).Value = 147
$ws.Cells.Item(153,3).Value = 1
$ws.Cells.Item(153,5).Value = 4
